# Daily attendance processing - 2026-01-29 23:41:47
# Reorders the "Recorded By" (column G) contributor list for rows whose
# value exactly matches one of the two known orderings, rotating the
# list right by one (last contributor moved to the front).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "backup@backdoor.com, System, system") {
        $cell.Value = "system, backup@backdoor.com, System"
    }
}
